# Refresh Universalis market-price snapshots + dependent Leve-profit figures
# across all eight crafter-job sheets (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR).
# Columns: H=currentAveragePrice, I=currentAveragePriceNQ, J=currentAveragePriceHQ,
#          K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 11: Gotta Bounce
$ws.Cells.Item(11, 8).Value = 259.4
$ws.Cells.Item(11, 9).Value = 259.4
$ws.Cells.Item(11, 11).Value = 259.4
$ws.Cells.Item(11, 13).Value = -119.4

# Row 33: Glazed and Confused
$ws.Cells.Item(33, 8).Value = 519.73914
$ws.Cells.Item(33, 9).Value = 497
$ws.Cells.Item(33, 10).Value = 562.375
$ws.Cells.Item(33, 11).Value = 497
$ws.Cells.Item(33, 12).Value = 562.375
$ws.Cells.Item(33, 13).Value = -268
$ws.Cells.Item(33, 14).Value = -1020.375

# Row 86: Filling in the Blanks
$ws.Cells.Item(86, 8).Value = 4442
$ws.Cells.Item(86, 9).Value = 3699
$ws.Cells.Item(86, 11).Value = 3699
$ws.Cells.Item(86, 13).Value = -2576

# Row 89: Ink into Antiquity (L)
$ws.Cells.Item(89, 8).Value = 4442
$ws.Cells.Item(89, 9).Value = 3699
$ws.Cells.Item(89, 11).Value = 18495
$ws.Cells.Item(89, 13).Value = -12879

# Row 137: Cutting Edge of Culinary Quality
$ws.Cells.Item(137, 8).Value = 23258752
$ws.Cells.Item(137, 9).Value = 50002216
$ws.Cells.Item(137, 11).Value = 150006648
$ws.Cells.Item(137, 13).Value = -150004098

# Row 138: All-night Crafting
$ws.Cells.Item(138, 8).Value = 4482.59
$ws.Cells.Item(138, 9).Value = 2522.0588
$ws.Cells.Item(138, 10).Value = 5997.5454
$ws.Cells.Item(138, 11).Value = 7566.176399999999
$ws.Cells.Item(138, 12).Value = 17992.6362
$ws.Cells.Item(138, 13).Value = -2426.176399999999
$ws.Cells.Item(138, 14).Value = -28272.6362

# Row 141: Remedy for Reason
$ws.Cells.Item(141, 8).Value = 1159.6786
$ws.Cells.Item(141, 9).Value = 979.6923
$ws.Cells.Item(141, 11).Value = 2939.0769
$ws.Cells.Item(141, 13).Value = 2240.9231

$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots
$ws.Cells.Item(2, 8).Value = 6464.1665
$ws.Cells.Item(2, 10).Value = 21341.4
$ws.Cells.Item(2, 12).Value = 21341.4
$ws.Cells.Item(2, 14).Value = -21567.4

# Row 18: Still the Best
$ws.Cells.Item(18, 8).Value = 256.5
$ws.Cells.Item(18, 10).Value = 256.5
$ws.Cells.Item(18, 12).Value = 256.5
$ws.Cells.Item(18, 14).Value = -900.5

# Row 25: Still Crazy After All These Years
$ws.Cells.Item(25, 8).Value = 813.4286
$ws.Cells.Item(25, 9).Value = 849.1667
$ws.Cells.Item(25, 10).Value = 599
$ws.Cells.Item(25, 11).Value = 849.1667
$ws.Cells.Item(25, 12).Value = 599
$ws.Cells.Item(25, 13).Value = -447.1667
$ws.Cells.Item(25, 14).Value = -1403

# Row 74: As the Bolt Flies
$ws.Cells.Item(74, 8).Value = 33336914
$ws.Cells.Item(74, 9).Value = 55558024
$ws.Cells.Item(74, 10).Value = 5247.25
$ws.Cells.Item(74, 11).Value = 55558024
$ws.Cells.Item(74, 12).Value = 5247.25
$ws.Cells.Item(74, 13).Value = -55557150
$ws.Cells.Item(74, 14).Value = -6995.25

# Row 77: Heavy Metal Banned (L)
$ws.Cells.Item(77, 8).Value = 33336914
$ws.Cells.Item(77, 9).Value = 55558024
$ws.Cells.Item(77, 10).Value = 5247.25
$ws.Cells.Item(77, 11).Value = 277790120
$ws.Cells.Item(77, 12).Value = 26236.25
$ws.Cells.Item(77, 13).Value = -277785752
$ws.Cells.Item(77, 14).Value = -34972.25

# Row 116: No Scope
$ws.Cells.Item(116, 8).Value = 6464.1665
$ws.Cells.Item(116, 10).Value = 21341.4
$ws.Cells.Item(116, 12).Value = 21341.4
$ws.Cells.Item(116, 14).Value = -25929.4

# Row 122: Haste for High Durium
$ws.Cells.Item(122, 8).Value = 55557656
$ws.Cells.Item(122, 9).Value = 1481.25
$ws.Cells.Item(122, 10).Value = 166670000
$ws.Cells.Item(122, 11).Value = 4443.75
$ws.Cells.Item(122, 12).Value = 500010000
$ws.Cells.Item(122, 13).Value = -1993.75
$ws.Cells.Item(122, 14).Value = -500014900

$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells
$ws.Cells.Item(3, 8).Value = 6464.1665
$ws.Cells.Item(3, 10).Value = 21341.4
$ws.Cells.Item(3, 12).Value = 21341.4
$ws.Cells.Item(3, 14).Value = -21569.4

# Row 5: Axe Me Anything
$ws.Cells.Item(5, 8).Value = 3055.75
$ws.Cells.Item(5, 9).Value = 899
$ws.Cells.Item(5, 10).Value = 3774.6667
$ws.Cells.Item(5, 11).Value = 899
$ws.Cells.Item(5, 12).Value = 3774.6667
$ws.Cells.Item(5, 13).Value = -786
$ws.Cells.Item(5, 14).Value = -4000.6667

# Row 134: Ruthenium Supremium
$ws.Cells.Item(134, 8).Value = 3425.1428
$ws.Cells.Item(134, 9).Value = 1338.375
$ws.Cells.Item(134, 10).Value = 10102.8
$ws.Cells.Item(134, 11).Value = 4015.125
$ws.Cells.Item(134, 12).Value = 30308.4
$ws.Cells.Item(134, 13).Value = -1480.125
$ws.Cells.Item(134, 14).Value = -35378.39999999999

# Row 138: Bladewinner
$ws.Cells.Item(138, 8).Value = 65499.25
$ws.Cells.Item(138, 10).Value = 65499.25
$ws.Cells.Item(138, 12).Value = 65499.25
$ws.Cells.Item(138, 14).Value = -75779.25

# Row 140: Ceremonial Teeth
$ws.Cells.Item(140, 8).Value = 57662.25
$ws.Cells.Item(140, 10).Value = 57662.25
$ws.Cells.Item(140, 12).Value = 57662.25
$ws.Cells.Item(140, 14).Value = -68022.25

# Row 141: Awl Dreams Come True
$ws.Cells.Item(141, 8).Value = 0
$ws.Cells.Item(141, 10).Value = 0
$ws.Cells.Item(141, 14).Value = 0

$ws = $wb.Worksheets.Item("CRP")
# Row 15: On the Move
$ws.Cells.Item(15, 8).Value = 7844
$ws.Cells.Item(15, 9).Value = 3067.2273
$ws.Cells.Item(15, 11).Value = 3067.2273
$ws.Cells.Item(15, 13).Value = -2897.2273

# Row 31: Wall Not Found
$ws.Cells.Item(31, 8).Value = 37238
$ws.Cells.Item(31, 9).Value = 2455.5715
$ws.Cells.Item(31, 10).Value = 46977.08
$ws.Cells.Item(31, 11).Value = 2455.5715
$ws.Cells.Item(31, 12).Value = 46977.08
$ws.Cells.Item(31, 13).Value = -2160.5715
$ws.Cells.Item(31, 14).Value = -47567.08

# Row 34: Armoires of the Rich and Famous
$ws.Cells.Item(34, 8).Value = 37238
$ws.Cells.Item(34, 9).Value = 2455.5715
$ws.Cells.Item(34, 10).Value = 46977.08
$ws.Cells.Item(34, 11).Value = 2455.5715
$ws.Cells.Item(34, 12).Value = 46977.08
$ws.Cells.Item(34, 13).Value = -2253.5715
$ws.Cells.Item(34, 14).Value = -47381.08

# Row 132: Hull Lotta Damage
$ws.Cells.Item(132, 8).Value = 2716
$ws.Cells.Item(132, 9).Value = 2236.6956
$ws.Cells.Item(132, 11).Value = 6710.0868
$ws.Cells.Item(132, 13).Value = -4180.0868

# Row 141: No Greater Treasure
$ws.Cells.Item(141, 8).Value = 241996
$ws.Cells.Item(141, 10).Value = 241996
$ws.Cells.Item(141, 12).Value = 241996
$ws.Cells.Item(141, 14).Value = -252356

$ws = $wb.Worksheets.Item("CUL")
# Row 4: In Hot Water
$ws.Cells.Item(4, 8).Value = 34690170
$ws.Cells.Item(4, 9).Value = 40645020
$ws.Cells.Item(4, 10).Value = 19532360
$ws.Cells.Item(4, 11).Value = 121935060
$ws.Cells.Item(4, 12).Value = 58597080
$ws.Cells.Item(4, 13).Value = -121934948
$ws.Cells.Item(4, 14).Value = -58597304

# Row 23: Sweet Smell of Success
$ws.Cells.Item(23, 8).Value = 798.1667
$ws.Cells.Item(23, 9).Value = 199.5
$ws.Cells.Item(23, 10).Value = 1097.5
$ws.Cells.Item(23, 11).Value = 598.5
$ws.Cells.Item(23, 12).Value = 3292.5
$ws.Cells.Item(23, 13).Value = -363.5
$ws.Cells.Item(23, 14).Value = -3762.5

# Row 131: The Mountain Steeped
$ws.Cells.Item(131, 8).Value = 6798561
$ws.Cells.Item(131, 9).Value = 13890780
$ws.Cells.Item(131, 10).Value = 5118825
$ws.Cells.Item(131, 11).Value = 41672340
$ws.Cells.Item(131, 12).Value = 15356475
$ws.Cells.Item(131, 13).Value = -41667300
$ws.Cells.Item(131, 14).Value = -15366555

$ws = $wb.Worksheets.Item("GSM")
# Row 39: One Man's Trash
$ws.Cells.Item(39, 8).Value = 14995
$ws.Cells.Item(39, 10).Value = 14995
$ws.Cells.Item(39, 12).Value = 14995
$ws.Cells.Item(39, 14).Value = -16059

# Row 132: On Board for Lar
$ws.Cells.Item(132, 8).Value = 3483.1194
$ws.Cells.Item(132, 9).Value = 3076.0544
$ws.Cells.Item(132, 10).Value = 5348.8335
$ws.Cells.Item(132, 11).Value = 9228.163199999999
$ws.Cells.Item(132, 12).Value = 16046.5005
$ws.Cells.Item(132, 13).Value = -6698.163199999999
$ws.Cells.Item(132, 14).Value = -21106.5005

$ws = $wb.Worksheets.Item("LTW")
# Row 16: Saddle Sore
$ws.Cells.Item(16, 8).Value = 655.9375
$ws.Cells.Item(16, 9).Value = 676.75
$ws.Cells.Item(16, 11).Value = 676.75
$ws.Cells.Item(16, 13).Value = -506.75

# Row 55: It's Not a Job, It's a Calling
$ws.Cells.Item(55, 8).Value = 3573873
$ws.Cells.Item(55, 9).Value = 6250677
$ws.Cells.Item(55, 10).Value = 4800.8335
$ws.Cells.Item(55, 11).Value = 6250677
$ws.Cells.Item(55, 12).Value = 4800.8335
$ws.Cells.Item(55, 13).Value = -6250504
$ws.Cells.Item(55, 14).Value = -5146.8335

# Row 81: I Need Your Glove Tonight
$ws.Cells.Item(81, 8).Value = 49000
$ws.Cells.Item(81, 10).Value = 49000
$ws.Cells.Item(81, 12).Value = 49000
$ws.Cells.Item(81, 14).Value = -50996

# Row 84: Halonic Drake Handlers (L)
$ws.Cells.Item(84, 8).Value = 49000
$ws.Cells.Item(84, 10).Value = 49000
$ws.Cells.Item(84, 12).Value = 147000
$ws.Cells.Item(84, 14).Value = -156984

# Row 87: Bar of the Bannermen
$ws.Cells.Item(87, 10).Value = 50000
$ws.Cells.Item(87, 12).Value = 50000
$ws.Cells.Item(87, 14).Value = -52246

# Row 90: Do My Little Turn on the Stonewalk (L)
$ws.Cells.Item(90, 10).Value = 50000
$ws.Cells.Item(90, 12).Value = 150000
$ws.Cells.Item(90, 14).Value = -161232

# Row 132: Tenets of Tanning
$ws.Cells.Item(132, 8).Value = 6705.5356
$ws.Cells.Item(132, 9).Value = 5485.6665
$ws.Cells.Item(132, 10).Value = 14024.75
$ws.Cells.Item(132, 11).Value = 16456.9995
$ws.Cells.Item(132, 12).Value = 42074.25
$ws.Cells.Item(132, 13).Value = -13926.9995
$ws.Cells.Item(132, 14).Value = -47134.25

# Row 136: Respect for Br'aax
$ws.Cells.Item(136, 8).Value = 6076.6665
$ws.Cells.Item(136, 9).Value = 3913.9333
$ws.Cells.Item(136, 10).Value = 16890.334
$ws.Cells.Item(136, 11).Value = 11741.7999
$ws.Cells.Item(136, 12).Value = 50671.00199999999
$ws.Cells.Item(136, 13).Value = -9191.7999
$ws.Cells.Item(136, 14).Value = -55771.00199999999

$ws = $wb.Worksheets.Item("WVR")
# Row 12: This Is Why You Can't Have Nice Things
$ws.Cells.Item(12, 8).Value = 0
$ws.Cells.Item(12, 10).Value = 0
$ws.Cells.Item(12, 14).Value = 0
